$wb = $excel.ActiveWorkbook

# ============================================================
# Sheet 1: rename existing sheet to "Sales vs PO" and restructure
# ============================================================
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Sales vs PO"

# Insert a new column before column C -> old C (PO_Requested_Qty) becomes D
$ws1.Columns.Item(3).Insert()

# New column C header text (cell already carries header style s=1 from the insert)
$ws1.Range("C1").Value = "Order Week"

$oldA = @(45551,45558,45565,45572,45579,45586,45593,45600,45607,45614,45621,45628,45635,45642,45649)
$newA = @(45557,45564,45571,45578,45585,45592,45599,45606,45613,45620,45627,45634,45641,45648,45655)

for ($i = 0; $i -lt $oldA.Length; $i++) {
    $row = $i + 2
    $ws1.Range("A$row").Value = $newA[$i]
    $ws1.Range("C$row").Value = $oldA[$i]
    $ws1.Range("D$row").Value = 0
}

# Give new column C data cells the same date style as column A (copy format only)
$ws1.Range("A2:A16").Copy($ws1.Range("C2:C16"))
for ($i = 0; $i -lt $oldA.Length; $i++) {
    $row = $i + 2
    $ws1.Range("C$row").Value = $oldA[$i]
}

# ============================================================
# Sheet 2: Weekly Growth  (copy sheet1's formatting/structure, then overwrite)
# ============================================================
$ws1.Copy($null, $ws1) | Out-Null
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "Weekly Growth"
$ws2.Cells.Clear()

# Bring header style (s=1) and date style (s=2) over from sheet1
$ws1.Range("A1:C1").Copy($ws2.Range("A1:C1"))
$ws1.Range("A2:A4").Copy($ws2.Range("A2:A4"))

$ws2.Range("A1").Value = "ds"
$ws2.Range("B1").Value = "PO_Requested_Qty"
$ws2.Range("C1").Value = "Growth%"

$ws2.Range("A2").Value = 45558
$ws2.Range("B2").Value = 992
$ws2.Range("C2").Value = 0

$ws2.Range("A3").Value = 45593
$ws2.Range("B3").Value = 16
$ws2.Range("C3").Value = -98.38709677419355

$ws2.Range("A4").Value = 45628
$ws2.Range("B4").Value = 672
$ws2.Range("C4").Value = 4100

# ============================================================
# Sheet 3: Volume Insights
# ============================================================
$ws1.Copy($null, $ws2) | Out-Null
$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "Volume Insights"
$ws3.Cells.Clear()

$ws1.Range("A1:D1").Copy($ws3.Range("A1:D1"))

$ws3.Range("A1").Value = "Total_PO_Quantity"
$ws3.Range("B1").Value = "Average_PO_Quantity"
$ws3.Range("C1").Value = "Max_PO_Quantity"
$ws3.Range("D1").Value = "Min_PO_Quantity"

$ws3.Range("A2").Value = 1680
$ws3.Range("B2").Value = 560
$ws3.Range("C2").Value = 992
$ws3.Range("D2").Value = 16

# ============================================================
# Sheet 4: Prediction Info
# ============================================================
$ws1.Copy($null, $ws3) | Out-Null
$ws4 = $wb.Worksheets.Item(4)
$ws4.Name = "Prediction Info"
$ws4.Cells.Clear()

$ws1.Range("A1").Copy($ws4.Range("A1"))

$ws4.Range("A1").Value = "Predicted_Next_Week_PO_Quantity"
$ws4.Range("A2").Value = 240.0000000000001

# ============================================================
# Restore the first sheet as the active sheet/tab
# ============================================================
$ws1.Activate()
